$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 142858620
$ws.Range("I40").Value = 1400
$ws.Range("K40").Value = 1400
$ws.Range("M40").Value = -1225
$ws.Range("H74").Value = 6359.1177
$ws.Range("I74").Value = 6748.6206
$ws.Range("J74").Value = 4100
$ws.Range("K74").Value = 6748.6206
$ws.Range("L74").Value = 4100
$ws.Range("M74").Value = -5812.6206
$ws.Range("N74").Value = -5972
$ws.Range("H77").Value = 6359.1177
$ws.Range("I77").Value = 6748.6206
$ws.Range("J77").Value = 4100
$ws.Range("K77").Value = 33743.103
$ws.Range("L77").Value = 20500
$ws.Range("M77").Value = -29063.103
$ws.Range("N77").Value = -29860

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2562.6667
$ws.Range("I2").Value = 1531.2
$ws.Range("J2").Value = 7720
$ws.Range("K2").Value = 1531.2
$ws.Range("L2").Value = 7720
$ws.Range("M2").Value = -1418.2
$ws.Range("N2").Value = -7946
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 3105175.8
$ws.Range("I32").Value = 3285.4429
$ws.Range("K32").Value = 3285.4429
$ws.Range("M32").Value = -2998.4429
$ws.Range("H116").Value = 2562.6667
$ws.Range("I116").Value = 1531.2
$ws.Range("J116").Value = 7720
$ws.Range("K116").Value = 1531.2
$ws.Range("L116").Value = 7720
$ws.Range("M116").Value = 762.8
$ws.Range("N116").Value = -12308
$ws.Range("H132").Value = 148725.62
$ws.Range("I132").Value = 157708.53
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 473125.59
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -470595.59
$ws.Range("N132").Value = -20057

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2562.6667
$ws.Range("I3").Value = 1531.2
$ws.Range("J3").Value = 7720
$ws.Range("K3").Value = 1531.2
$ws.Range("L3").Value = 7720
$ws.Range("M3").Value = -1417.2
$ws.Range("N3").Value = -7948
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H13").Value = 38333.332
$ws.Range("J13").Value = 38333.332
$ws.Range("L13").Value = 38333.332
$ws.Range("N13").Value = -38669.332
$ws.Range("H50").Value = 14180
$ws.Range("J50").Value = 14180
$ws.Range("L50").Value = 14180
$ws.Range("N50").Value = -15328
$ws.Range("H86").Value = 2625.75
$ws.Range("I86").Value = 2851
$ws.Range("J86").Value = 1950
$ws.Range("K86").Value = 2851
$ws.Range("L86").Value = 1950
$ws.Range("M86").Value = -1728
$ws.Range("N86").Value = -4196
$ws.Range("H89").Value = 2625.75
$ws.Range("I89").Value = 2851
$ws.Range("J89").Value = 1950
$ws.Range("K89").Value = 14255
$ws.Range("L89").Value = 9750
$ws.Range("M89").Value = -8639
$ws.Range("N89").Value = -20982
$ws.Range("H94").Value = 620.8570999999999
$ws.Range("I94").Value = 620.8570999999999
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 620.8570999999999
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -169.8570999999999
$ws.Range("N94").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1838.7667
$ws.Range("I31").Value = 1605.8235
$ws.Range("J31").Value = 2143.3845
$ws.Range("K31").Value = 1605.8235
$ws.Range("L31").Value = 2143.3845
$ws.Range("M31").Value = -1310.8235
$ws.Range("N31").Value = -2733.3845
$ws.Range("H34").Value = 1838.7667
$ws.Range("I34").Value = 1605.8235
$ws.Range("J34").Value = 2143.3845
$ws.Range("K34").Value = 1605.8235
$ws.Range("L34").Value = 2143.3845
$ws.Range("M34").Value = -1403.8235
$ws.Range("N34").Value = -2547.3845
$ws.Range("H59").Value = 36368.8
$ws.Range("J59").Value = 36368.8
$ws.Range("L59").Value = 36368.8
$ws.Range("N59").Value = -38658.8
$ws.Range("H62").Value = 2957.1428
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 3116.6667
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 3116.6667
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4364.6667
$ws.Range("H65").Value = 2957.1428
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 3116.6667
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 15583.3335
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -21823.3335
$ws.Range("H68").Value = 18692.5
$ws.Range("J68").Value = 18692.5
$ws.Range("L68").Value = 18692.5
$ws.Range("N68").Value = -20190.5
$ws.Range("H70").Value = 23495
$ws.Range("J70").Value = 23495
$ws.Range("L70").Value = 23495
$ws.Range("N70").Value = -24125
$ws.Range("H71").Value = 18692.5
$ws.Range("J71").Value = 18692.5
$ws.Range("L71").Value = 56077.5
$ws.Range("N71").Value = -63565.5
$ws.Range("H73").Value = 23495
$ws.Range("J73").Value = 23495
$ws.Range("L73").Value = 23495
$ws.Range("N73").Value = -25679
$ws.Range("H74").Value = 25877.111
$ws.Range("J74").Value = 25877.111
$ws.Range("L74").Value = 25877.111
$ws.Range("N74").Value = -27625.111
$ws.Range("H75").Value = 37000
$ws.Range("J75").Value = 37000
$ws.Range("L75").Value = 37000
$ws.Range("N75").Value = -38996
$ws.Range("H77").Value = 25877.111
$ws.Range("J77").Value = 25877.111
$ws.Range("L77").Value = 77631.333
$ws.Range("N77").Value = -86367.333
$ws.Range("H78").Value = 37000
$ws.Range("J78").Value = 37000
$ws.Range("L78").Value = 111000
$ws.Range("N78").Value = -120984
$ws.Range("H82").Value = 44000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 44000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 44000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -44722
$ws.Range("H85").Value = 44000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 44000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 44000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -46496
$ws.Range("H88").Value = 31722
$ws.Range("J88").Value = 31722
$ws.Range("L88").Value = 31722
$ws.Range("N88").Value = -32534
$ws.Range("H91").Value = 31722
$ws.Range("J91").Value = 31722
$ws.Range("L91").Value = 31722
$ws.Range("N91").Value = -34530
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H132").Value = 3186.8572
$ws.Range("I132").Value = 2634.2856
$ws.Range("J132").Value = 4844.5713
$ws.Range("K132").Value = 7902.8568
$ws.Range("L132").Value = 14533.7139
$ws.Range("M132").Value = -5372.8568
$ws.Range("N132").Value = -19593.7139

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2979081.5
$ws.Range("I80").Value = 2697.6667
$ws.Range("J80").Value = 3971209.2
$ws.Range("K80").Value = 2697.6667
$ws.Range("L80").Value = 3971209.2
$ws.Range("M80").Value = -1699.6667
$ws.Range("N80").Value = -3973205.2
$ws.Range("H83").Value = 2979081.5
$ws.Range("I83").Value = 2697.6667
$ws.Range("J83").Value = 3971209.2
$ws.Range("K83").Value = 13488.3335
$ws.Range("L83").Value = 19856046
$ws.Range("M83").Value = -8496.333500000001
$ws.Range("N83").Value = -19866030
$ws.Range("H113").Value = 1943.421
$ws.Range("I113").Value = 1813.0625
$ws.Range("J113").Value = 2638.6667
$ws.Range("K113").Value = 1813.0625
$ws.Range("L113").Value = 2638.6667
$ws.Range("M113").Value = 356.9375
$ws.Range("N113").Value = -6978.6667
$ws.Range("H122").Value = 3537.0588
$ws.Range("I122").Value = 5178.3335
$ws.Range("J122").Value = 2641.818
$ws.Range("K122").Value = 15535.0005
$ws.Range("L122").Value = 7925.454000000001
$ws.Range("M122").Value = -13085.0005
$ws.Range("N122").Value = -12825.454
$ws.Range("H132").Value = 2711.5925
$ws.Range("I132").Value = 2173.6
$ws.Range("J132").Value = 3028.0588
$ws.Range("K132").Value = 6520.799999999999
$ws.Range("L132").Value = 9084.1764
$ws.Range("M132").Value = -3990.799999999999
$ws.Range("N132").Value = -14144.1764

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1438.5385
$ws.Range("I46").Value = 1154.6364
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 1154.6364
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -966.6364000000001
$ws.Range("N46").Value = -3376
$ws.Range("H122").Value = 2724.5
$ws.Range("I122").Value = 2213.8462
$ws.Range("J122").Value = 3073.8948
$ws.Range("K122").Value = 6641.5386
$ws.Range("L122").Value = 9221.6844
$ws.Range("M122").Value = -4191.5386
$ws.Range("N122").Value = -14121.6844
$ws.Range("H132").Value = 2696.9285
$ws.Range("I132").Value = 2562.875
$ws.Range("J132").Value = 2875.6667
$ws.Range("K132").Value = 7688.625
$ws.Range("L132").Value = 8627.000100000001
$ws.Range("M132").Value = -5158.625
$ws.Range("N132").Value = -13687.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H113").Value = 44572.695
$ws.Range("I113").Value = 67013.60000000001
$ws.Range("J113").Value = 2496
$ws.Range("K113").Value = 201040.8
$ws.Range("L113").Value = 7488
$ws.Range("M113").Value = -198870.8
$ws.Range("N113").Value = -11828
$ws.Range("H122").Value = 3536.7144
$ws.Range("I122").Value = 2475.5
$ws.Range("J122").Value = 4951.6665
$ws.Range("K122").Value = 7426.5
$ws.Range("L122").Value = 14854.9995
$ws.Range("M122").Value = -4976.5
$ws.Range("N122").Value = -19754.9995
$ws.Range("H132").Value = 1739.841
$ws.Range("I132").Value = 1583.5676
$ws.Range("J132").Value = 2565.8572
$ws.Range("K132").Value = 4750.7028
$ws.Range("L132").Value = 7697.571599999999
$ws.Range("M132").Value = -2220.7028
$ws.Range("N132").Value = -12757.5716
